# Fruta / hortaliza, semanal
# Update the Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) columns for the weekly
# refresh of this market/product subset. Row 10 is untouched; rows 2-9 and
# 11-13 are refreshed with the new weekly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 45096; M = 30 },
    @{ Row = 3;  D = 45055; M = 50; N = 15000; O = 15000; P = 15000; S = 833 },
    @{ Row = 4;  D = 45069; M = 60 },
    @{ Row = 5;  D = 45085; M = 30; N = 19000; P = 19000; S = 1056 },
    @{ Row = 6;  D = 45111; M = 20; N = 20000; O = 20000; P = 20000; S = 1111 },
    @{ Row = 7;  D = 45083; M = 50 },
    @{ Row = 8;  D = 45072; M = 30 },
    @{ Row = 9;  D = 45061; M = 40; N = 15000; O = 15000; P = 15000; S = 833 },
    @{ Row = 11; D = 45076; M = 20 },
    @{ Row = 12; D = 45092; N = 18000; O = 19000; P = 18667; S = 1037 },
    @{ Row = 13; D = 45084; M = 50; N = 18000; O = 19000; P = 18500; S = 1028 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    if ($u.ContainsKey("M")) { $ws.Range("M$r").Value = $u.M }
    if ($u.ContainsKey("N")) { $ws.Range("N$r").Value = $u.N }
    if ($u.ContainsKey("O")) { $ws.Range("O$r").Value = $u.O }
    if ($u.ContainsKey("P")) { $ws.Range("P$r").Value = $u.P }
    if ($u.ContainsKey("S")) { $ws.Range("S$r").Value = $u.S }
}
